# Auto-generated Excel COM-interop script applying the Anima_Profits.xlsx data refresh
# (scheduled runner update: recalculated currentAveragePrice / LevePrice* / LeveProfit* columns)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H7").Value = 80006
$ws.Range("J7").Value = 80006
$ws.Range("L7").Value = 80006
$ws.Range("N7").Value = -80230
$ws.Range("H14").Value = 80006
$ws.Range("J14").Value = 80006
$ws.Range("L14").Value = 80006
$ws.Range("N14").Value = -80388
$ws.Range("H62").Value = 2728.3333
$ws.Range("I62").Value = 1656.6666
$ws.Range("K62").Value = 1656.6666
$ws.Range("M62").Value = -1032.6666
$ws.Range("H65").Value = 2728.3333
$ws.Range("I65").Value = 1656.6666
$ws.Range("K65").Value = 8283.333000000001
$ws.Range("M65").Value = -5163.333000000001
$ws.Range("H94").Value = 2096.6667
$ws.Range("I94").Value = 2096.6667
$ws.Range("K94").Value = 2096.6667
$ws.Range("M94").Value = -1645.6667
$ws.Range("H138").Value = 207138.8
$ws.Range("I138").Value = 2725.6667
$ws.Range("J138").Value = 287828.2
$ws.Range("K138").Value = 8177.000100000001
$ws.Range("L138").Value = 863484.6000000001
$ws.Range("M138").Value = -3037.000100000001
$ws.Range("N138").Value = -873764.6000000001

$ws = $wb.Worksheets("ARM")
$ws.Range("H38").Value = 61500
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 61500
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 61500
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -62434
$ws.Range("H74").Value = 2293.6
$ws.Range("I74").Value = 1992
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 1992
$ws.Range("L74").Value = 3500
$ws.Range("M74").Value = -1118
$ws.Range("N74").Value = -5248
$ws.Range("H77").Value = 2293.6
$ws.Range("I77").Value = 1992
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 9960
$ws.Range("L77").Value = 17500
$ws.Range("M77").Value = -5592
$ws.Range("N77").Value = -26236
$ws.Range("H122").Value = 2634.2666
$ws.Range("I122").Value = 2348.6667
$ws.Range("K122").Value = 7046.000100000001
$ws.Range("M122").Value = -4596.000100000001
$ws.Range("H132").Value = 5991.4062
$ws.Range("I132").Value = 6206.125
$ws.Range("J132").Value = 5776.6875
$ws.Range("K132").Value = 18618.375
$ws.Range("L132").Value = 17330.0625
$ws.Range("M132").Value = -16088.375
$ws.Range("N132").Value = -22390.0625

$ws = $wb.Worksheets("BSM")
$ws.Range("H107").Value = 1904.1111
$ws.Range("I107").Value = 970.3333
$ws.Range("J107").Value = 2371
$ws.Range("K107").Value = 970.3333
$ws.Range("L107").Value = 2371
$ws.Range("M107").Value = 949.6667
$ws.Range("N107").Value = -6211

$ws = $wb.Worksheets("CRP")
$ws.Range("H3").Value = 37125.75
$ws.Range("J3").Value = 48501
$ws.Range("L3").Value = 48501
$ws.Range("N3").Value = -48727
$ws.Range("H107").Value = 2604733.2
$ws.Range("I107").Value = 3676926.2
$ws.Range("J107").Value = 835.7143
$ws.Range("K107").Value = 3676926.2
$ws.Range("L107").Value = 835.7143
$ws.Range("M107").Value = -3675006.2
$ws.Range("N107").Value = -4675.7143
$ws.Range("H111").Value = 38793.57
$ws.Range("J111").Value = 38793.57
$ws.Range("L111").Value = 38793.57
$ws.Range("N111").Value = -46973.57
$ws.Range("H122").Value = 1553.6428
$ws.Range("I122").Value = 863.8182
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 2591.4546
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -141.4546
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets("CUL")
$ws.Range("H2").Value = 64.07143000000001
$ws.Range("J2").Value = 118.333336
$ws.Range("L2").Value = 710.000016
$ws.Range("N2").Value = -936.000016
$ws.Range("H92").Value = 530.7368
$ws.Range("I92").Value = 424.4
$ws.Range("J92").Value = 568.7143
$ws.Range("K92").Value = 1273.2
$ws.Range("L92").Value = 1706.1429
$ws.Range("M92").Value = -25.19999999999982
$ws.Range("N92").Value = -4202.1429
$ws.Range("H132").Value = 3151.55
$ws.Range("J132").Value = 3811.9583
$ws.Range("L132").Value = 34307.6247
$ws.Range("N132").Value = -39367.6247
$ws.Range("H136").Value = 2128.7896
$ws.Range("I136").Value = 1249.6666
$ws.Range("J136").Value = 2920
$ws.Range("K136").Value = 3748.9998
$ws.Range("L136").Value = 8760
$ws.Range("M136").Value = 1351.0002
$ws.Range("N136").Value = -18960
$ws.Range("H137").Value = 11913961
$ws.Range("I137").Value = 13899129
$ws.Range("J137").Value = 2950
$ws.Range("K137").Value = 41697387
$ws.Range("L137").Value = 8850
$ws.Range("M137").Value = -41692287
$ws.Range("N137").Value = -19050

$ws = $wb.Worksheets("GSM")
$ws.Range("H41").Value = 2959.2
$ws.Range("I41").Value = 3249
$ws.Range("J41").Value = 1800
$ws.Range("K41").Value = 3249
$ws.Range("L41").Value = 1800
$ws.Range("M41").Value = -2894
$ws.Range("N41").Value = -2510
$ws.Range("H122").Value = 5448.9165
$ws.Range("I122").Value = 4459
$ws.Range("J122").Value = 5709.421
$ws.Range("K122").Value = 13377
$ws.Range("L122").Value = 17128.263
$ws.Range("M122").Value = -10927
$ws.Range("N122").Value = -22028.263
$ws.Range("H126").Value = 2078
$ws.Range("I126").Value = 1351
$ws.Range("J126").Value = 2659.6
$ws.Range("K126").Value = 4053
$ws.Range("L126").Value = 7978.799999999999
$ws.Range("M126").Value = -1583
$ws.Range("N126").Value = -12918.8

$ws = $wb.Worksheets("LTW")
$ws.Range("H2").Value = 24857.715
$ws.Range("H18").Value = 70006
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H100").Value = 2599.1428
$ws.Range("I100").Value = 2536.182
$ws.Range("J100").Value = 2830
$ws.Range("K100").Value = 2536.182
$ws.Range("L100").Value = 2830
$ws.Range("M100").Value = -1995.182
$ws.Range("N100").Value = -3912
$ws.Range("H122").Value = 2002
$ws.Range("I122").Value = 1004
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 3012
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -562
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 4434.069
$ws.Range("I132").Value = 3872.111
$ws.Range("J132").Value = 5353.636
$ws.Range("K132").Value = 11616.333
$ws.Range("L132").Value = 16060.908
$ws.Range("M132").Value = -9086.332999999999
$ws.Range("N132").Value = -21120.908

$ws = $wb.Worksheets("WVR")
$ws.Range("H3").Value = 25486.857
$ws.Range("I3").Value = 6600
$ws.Range("K3").Value = 6600
$ws.Range("M3").Value = -6486
$ws.Range("H20").Value = 11250
$ws.Range("J20").Value = 11250
$ws.Range("L20").Value = 11250
$ws.Range("N20").Value = -11730
$ws.Range("H32").Value = 333337340
$ws.Range("I32").Value = 6000
$ws.Range("J32").Value = 1000000000
$ws.Range("K32").Value = 6000
$ws.Range("L32").Value = 1000000000
$ws.Range("M32").Value = -5683
$ws.Range("N32").Value = -1000000634
$ws.Range("H122").Value = 3600.9092
$ws.Range("I122").Value = 1200
$ws.Range("K122").Value = 3600
$ws.Range("M122").Value = -1150
